$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.547304093837738
$ws.Range("B1").Value = 3.407247304916382
$ws.Range("C1").Value = 4.549694061279297
$ws.Range("D1").Value = 1.864233613014221
$ws.Range("E1").Value = 0.7921881675720215
